$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dados")

# Remove the obsolete first test row (old row 2); rows 3-11 shift up to 2-10
$ws.Rows(2).Delete()

# Re-anchor the "Altera" hyperlink AutoShapes (they sit one row below the
# deleted "0000DEPRE" row and must move up with it)
for ($i = 1; $i -le $ws.Shapes.Count; $i++) {
    $shp = $ws.Shapes.Item($i)
    $shp.Top = $shp.Top - $ws.Rows(1).RowHeight
}

# centroCusto column now uses the single test cost-center for every row
$ws.Range("A2:A10").Value = "_TST1_CC10"

# cpf column switches from text placeholders to real numeric CPF-like values
$ws.Range("B2").Value = 123454
$ws.Range("B3").Value = 123455
$ws.Range("B4").Value = 123456
$ws.Range("B5").Value = 123457
$ws.Range("B6").Value = 123458
$ws.Range("B7").Value = 123459
$ws.Range("B8").Value = 123460
$ws.Range("B9").Value = 123461
$ws.Range("B10").Value = 123462

# salarioBase for the last row now carries decimals
$ws.Range("F10").Value = 1000.34

$ws.Range("C8").Select()
